$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range('A2').Value = 'Última actualización: 04:30:03'
$ws1.Range('A3').Value = 'Total filas: 24'
$ws1.Range('A6').Value = '04:01:06'
$ws1.Range('B6').Value = '04:02'
$ws1.Range('C6').Value = '81_EL PELIGRO'
$ws1.Range('D6').Value = 1
$ws1.Range('A7').Value = '00:46:06'
$ws1.Range('B7').Value = '01:12'
$ws1.Range('C7').Value = '215_ALUAR'
$ws1.Range('D7').Value = 26
$ws1.Range('A8').Value = '04:01:06'
$ws1.Range('B8').Value = '04:47'
$ws1.Range('C8').Value = '215_EL PELIGRO'
$ws1.Range('D8').Value = 46
$ws1.Range('A9').Value = '03:46:12'
$ws1.Range('B9').Value = '04:46'
$ws1.Range('C9').Value = '215A_EL PATO'
$ws1.Range('D9').Value = 60
$ws1.Range('A10').Value = '01:55:38'
$ws1.Range('B10').Value = '03:02'
$ws1.Range('C10').Value = '15_ABASTO'
$ws1.Range('D10').Value = 67
$ws1.Range('A11').Value = '04:01:06'
$ws1.Range('B11').Value = '05:12'
$ws1.Range('C11').Value = '17_ROMERO'
$ws1.Range('D11').Value = 71
$ws1.Range('A12').Value = '00:46:06'
$ws1.Range('B12').Value = '01:58'
$ws1.Range('C12').Value = '14_ABASTO'
$ws1.Range('D12').Value = 72
$ws1.Range('A13').Value = '04:30:03'
$ws1.Range('B13').Value = '05:44'
$ws1.Range('C13').Value = '14_ABASTO'
$ws1.Range('D13').Value = 74
$ws1.Range('A14').Value = '03:46:12'
$ws1.Range('B14').Value = '05:16'
$ws1.Range('C14').Value = '17_ROMERO'
$ws1.Range('D14').Value = 90
$ws1.Range('A15').Value = '04:01:06'
$ws1.Range('B15').Value = '05:32'
$ws1.Range('C15').Value = '81_EL PELIGRO'
$ws1.Range('D15').Value = 91
$ws1.Range('A16').Value = '04:30:03'
$ws1.Range('B16').Value = '06:01'
$ws1.Range('C16').Value = '16_SANTA ANA'
$ws1.Range('D16').Value = 91
$ws1.Range('A17').Value = '02:29:13'
$ws1.Range('B17').Value = '04:01'
$ws1.Range('C17').Value = '81_EL PELIGRO'
$ws1.Range('D17').Value = 92
$ws1.Range('A18').Value = '04:30:03'
$ws1.Range('B18').Value = '06:04'
$ws1.Range('C18').Value = '10_OLMOS'
$ws1.Range('D18').Value = 94
$ws1.Range('A19').Value = '03:46:12'
$ws1.Range('B19').Value = '05:22'
$ws1.Range('C19').Value = '23_HERNANDEZ'
$ws1.Range('D19').Value = 96
$ws1.Range('A20').Value = '01:22:42'
$ws1.Range('B20').Value = '02:58'
$ws1.Range('C20').Value = '215_ALUAR'
$ws1.Range('D20').Value = 96
$ws1.Range('A21').Value = '04:30:03'
$ws1.Range('B21').Value = '06:11'
$ws1.Range('C21').Value = '215A_EL PATO'
$ws1.Range('D21').Value = 101
$ws1.Range('A22').Value = '04:01:06'
$ws1.Range('B22').Value = '05:45'
$ws1.Range('C22').Value = '14_ABASTO'
$ws1.Range('D22').Value = 104
$ws1.Range('A23').Value = '03:46:12'
$ws1.Range('B23').Value = '05:35'
$ws1.Range('C23').Value = '215B_EL PATO'
$ws1.Range('D23').Value = 109
$ws1.Range('A24').Value = '04:01:06'
$ws1.Range('B24').Value = '05:52'
$ws1.Range('C24').Value = '17_ROMERO'
$ws1.Range('D24').Value = 111
$ws1.Range('A25').Value = '01:55:38'
$ws1.Range('B25').Value = '03:48'
$ws1.Range('C25').Value = '14_ABASTO'
$ws1.Range('D25').Value = 113
$ws1.Range('A26').Value = '03:00:53'
$ws1.Range('B26').Value = '04:53'
$ws1.Range('C26').Value = '11_ETCHEVERRY'
$ws1.Range('D26').Value = 113
$ws1.Range('A27').Value = '04:30:03'
$ws1.Range('B27').Value = '06:24'
$ws1.Range('C27').Value = '11_ETCHEVERRY'
$ws1.Range('D27').Value = 114
$ws1.Range('A28').Value = '04:30:03'
$ws1.Range('B28').Value = '06:27'
$ws1.Range('C28').Value = '23_HERNANDEZ'
$ws1.Range('D28').Value = 117
$ws1.Range('A29').Value = '02:47:42'
$ws1.Range('B29').Value = '04:45'
$ws1.Range('C29').Value = '215A_EL PATO'
$ws1.Range('D29').Value = 118

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range('A2').Value = 'Última actualización: 04:30:03'
$ws2.Range('A3').Value = 'Total filas: 7'
$ws2.Range('A6').Value = '00:46:06'
$ws2.Range('B6').Value = '01:12'
$ws2.Range('C6').Value = '215_ALUAR'
$ws2.Range('D6').Value = 26
$ws2.Range('A7').Value = '04:01:06'
$ws2.Range('B7').Value = '04:47'
$ws2.Range('C7').Value = '215_EL PELIGRO'
$ws2.Range('D7').Value = 46
$ws2.Range('A8').Value = '03:46:12'
$ws2.Range('B8').Value = '04:46'
$ws2.Range('C8').Value = '215A_EL PATO'
$ws2.Range('D8').Value = 60
$ws2.Range('A9').Value = '01:22:42'
$ws2.Range('B9').Value = '02:58'
$ws2.Range('C9').Value = '215_ALUAR'
$ws2.Range('D9').Value = 96
$ws2.Range('A10').Value = '04:30:03'
$ws2.Range('B10').Value = '06:11'
$ws2.Range('C10').Value = '215A_EL PATO'
$ws2.Range('D10').Value = 101
$ws2.Range('A11').Value = '03:46:12'
$ws2.Range('B11').Value = '05:35'
$ws2.Range('C11').Value = '215B_EL PATO'
$ws2.Range('D11').Value = 109
$ws2.Range('A12').Value = '02:47:42'
$ws2.Range('B12').Value = '04:45'
$ws2.Range('C12').Value = '215A_EL PATO'
$ws2.Range('D12').Value = 118

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range('A2').Value = 'Última actualización: 04:30:03'
$ws3.Range('A3').Value = 'Total filas: 24'
$ws3.Range('A6').Value = '04:01:06'
$ws3.Range('B6').Value = '04:02'
$ws3.Range('C6').Value = '81_EL PELIGRO'
$ws3.Range('D6').Value = 1
$ws3.Range('A7').Value = '00:46:06'
$ws3.Range('B7').Value = '01:12'
$ws3.Range('C7').Value = '215_ALUAR'
$ws3.Range('D7').Value = 26
$ws3.Range('A8').Value = '04:01:06'
$ws3.Range('B8').Value = '04:47'
$ws3.Range('C8').Value = '215_EL PELIGRO'
$ws3.Range('D8').Value = 46
$ws3.Range('A9').Value = '03:46:12'
$ws3.Range('B9').Value = '04:46'
$ws3.Range('C9').Value = '215A_EL PATO'
$ws3.Range('D9').Value = 60
$ws3.Range('A10').Value = '01:55:38'
$ws3.Range('B10').Value = '03:02'
$ws3.Range('C10').Value = '15_ABASTO'
$ws3.Range('D10').Value = 67
$ws3.Range('A11').Value = '04:01:06'
$ws3.Range('B11').Value = '05:12'
$ws3.Range('C11').Value = '17_ROMERO'
$ws3.Range('D11').Value = 71
$ws3.Range('A12').Value = '00:46:06'
$ws3.Range('B12').Value = '01:58'
$ws3.Range('C12').Value = '14_ABASTO'
$ws3.Range('D12').Value = 72
$ws3.Range('A13').Value = '04:30:03'
$ws3.Range('B13').Value = '05:44'
$ws3.Range('C13').Value = '14_ABASTO'
$ws3.Range('D13').Value = 74
$ws3.Range('A14').Value = '03:46:12'
$ws3.Range('B14').Value = '05:16'
$ws3.Range('C14').Value = '17_ROMERO'
$ws3.Range('D14').Value = 90
$ws3.Range('A15').Value = '04:01:06'
$ws3.Range('B15').Value = '05:32'
$ws3.Range('C15').Value = '81_EL PELIGRO'
$ws3.Range('D15').Value = 91
$ws3.Range('A16').Value = '04:30:03'
$ws3.Range('B16').Value = '06:01'
$ws3.Range('C16').Value = '16_SANTA ANA'
$ws3.Range('D16').Value = 91
$ws3.Range('A17').Value = '02:29:13'
$ws3.Range('B17').Value = '04:01'
$ws3.Range('C17').Value = '81_EL PELIGRO'
$ws3.Range('D17').Value = 92
$ws3.Range('A18').Value = '04:30:03'
$ws3.Range('B18').Value = '06:04'
$ws3.Range('C18').Value = '10_OLMOS'
$ws3.Range('D18').Value = 94
$ws3.Range('A19').Value = '03:46:12'
$ws3.Range('B19').Value = '05:22'
$ws3.Range('C19').Value = '23_HERNANDEZ'
$ws3.Range('D19').Value = 96
$ws3.Range('A20').Value = '01:22:42'
$ws3.Range('B20').Value = '02:58'
$ws3.Range('C20').Value = '215_ALUAR'
$ws3.Range('D20').Value = 96
$ws3.Range('A21').Value = '04:30:03'
$ws3.Range('B21').Value = '06:11'
$ws3.Range('C21').Value = '215A_EL PATO'
$ws3.Range('D21').Value = 101
$ws3.Range('A22').Value = '04:01:06'
$ws3.Range('B22').Value = '05:45'
$ws3.Range('C22').Value = '14_ABASTO'
$ws3.Range('D22').Value = 104
$ws3.Range('A23').Value = '03:46:12'
$ws3.Range('B23').Value = '05:35'
$ws3.Range('C23').Value = '215B_EL PATO'
$ws3.Range('D23').Value = 109
$ws3.Range('A24').Value = '04:01:06'
$ws3.Range('B24').Value = '05:52'
$ws3.Range('C24').Value = '17_ROMERO'
$ws3.Range('D24').Value = 111
$ws3.Range('A25').Value = '01:55:38'
$ws3.Range('B25').Value = '03:48'
$ws3.Range('C25').Value = '14_ABASTO'
$ws3.Range('D25').Value = 113
$ws3.Range('A26').Value = '03:00:53'
$ws3.Range('B26').Value = '04:53'
$ws3.Range('C26').Value = '11_ETCHEVERRY'
$ws3.Range('D26').Value = 113
$ws3.Range('A27').Value = '04:30:03'
$ws3.Range('B27').Value = '06:24'
$ws3.Range('C27').Value = '11_ETCHEVERRY'
$ws3.Range('D27').Value = 114
$ws3.Range('A28').Value = '04:30:03'
$ws3.Range('B28').Value = '06:27'
$ws3.Range('C28').Value = '23_HERNANDEZ'
$ws3.Range('D28').Value = 117
$ws3.Range('A29').Value = '02:47:42'
$ws3.Range('B29').Value = '04:45'
$ws3.Range('C29').Value = '215A_EL PATO'
$ws3.Range('D29').Value = 118
